$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the confidence interval columns (G, H)
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# Confidence interval lower/upper bounds per row (rows without SE/df/t/p also
# have no CI values, matching the existing blank C:F cells on rows 12-14)
$ciData = @(
    @{ Row = 2;  Lower = -0.357047293664637;  Upper = 0.200392379310378 },
    @{ Row = 3;  Lower = 0.0508901659229239;  Upper = 0.15375273201626 },
    @{ Row = 4;  Lower = 0.0420926922049491;  Upper = 0.178125059615871 },
    @{ Row = 5;  Lower = 0.22292434843458;    Upper = 0.673513771927432 },
    @{ Row = 6;  Lower = 0.146656548956796;   Upper = 0.620611444827367 },
    @{ Row = 7;  Lower = 0.197516694468969;   Upper = 0.596750488874806 },
    @{ Row = 8;  Lower = 0.163371631332728;   Upper = 0.691345220424335 },
    @{ Row = 9;  Lower = 0.047527959665029;   Upper = 0.164902365214973 },
    @{ Row = 10; Lower = -0.207690556477201;  Upper = -0.152244880966598 },
    @{ Row = 11; Lower = -0.0354981485935685; Upper = 0.0199232947119327 },
    @{ Row = 15; Lower = -0.145659276928672;  Upper = -0.0511961940696945 },
    @{ Row = 16; Lower = -0.565337887505915;  Upper = -0.198704461631217 },
    @{ Row = 17; Lower = -0.712458927372385;  Upper = -0.248564256278552 },
    @{ Row = 18; Lower = -0.30750546403252;   Upper = -0.101780331845848 },
    @{ Row = 19; Lower = -1.19350097774961;   Upper = -0.395033388937938 },
    @{ Row = 20; Lower = -1.34702754385486;   Upper = -0.44584869686916 },
    @{ Row = 21; Lower = -0.0199232947119327; Upper = 0.0354981485935685 },
    @{ Row = 22; Lower = 0.0950559193300581;  Upper = 0.329804730429946 },
    @{ Row = 23; Lower = -0.0773269892731749; Upper = 0.137776657686462 },
    @{ Row = 24; Lower = 0.368934363551065;   Upper = 1.28004967154977 },
    @{ Row = 25; Lower = -0.152325749286223;  Upper = 0.0231556227083734 },
    @{ Row = 26; Lower = 0.37775870061599;    Upper = 1.28594741353018 }
)

foreach ($item in $ciData) {
    $ws.Cells.Item($item.Row, 7).Value = $item.Lower
    $ws.Cells.Item($item.Row, 8).Value = $item.Upper
}

# Rows 12-14 mirror the existing blank C:F cells - leave G/H blank but
# touch them so the cells materialize in the sheet (matches <c r="G12"/> etc.)
# Touching a formatting property (without actually changing the applied
# style) is enough to force the engine to emit the empty cell.
foreach ($row in 12..14) {
    $ws.Cells.Item($row, 7).Borders.LineStyle = -4142
    $ws.Cells.Item($row, 8).Borders.LineStyle = -4142
}
